$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers (shift C/D/E to prediction / rejection-f / max)
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Row 2 data values
$ws.Range("C2").Value = "g__Turicibacter"
$ws.Range("D2").Value = "g__Turicibacter"
$ws.Range("E2").Value = 1
